$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.106.62"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "2.291.97"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'317.60"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'102.74"
$ws.Range("E6").Value = "  -5.55%  "
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.603"
$ws.Range("E9").Value = "  -3.69%  "
$ws.Range("D10").Value = "'39.12"
$ws.Range("E10").Value = "  -6.90%  "
$ws.Range("D11").Value = "'0.0904"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "'0.962"
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("D15").Value = "'15.25"
$ws.Range("E15").Value = "  -5.85%  "
$ws.Range("D16").Value = "2.636.46"
$ws.Range("E16").Value = "  -3.37%  "
$ws.Range("D17").Value = "2.285.87"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "42.039.51"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("D19").Value = "'7.47"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "'3.67"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "'283.17"
$ws.Range("E22").Value = "  +9.60%  "
$ws.Range("D23").Value = "'73.32"
$ws.Range("E23").Value = "  -3.99%  "
$ws.Range("D24").Value = "'10.08"
$ws.Range("E24").Value = "  +6.80%  "
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "'10.78"
$ws.Range("E27").Value = "  -6.22%  "
$ws.Range("D28").Value = "'2.33"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("D29").Value = "'22.92"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "'35.63"
$ws.Range("E30").Value = "  -3.68%  "
$ws.Range("D31").Value = "'163.55"
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").Value = "'5.84"
$ws.Range("E33").Value = "  -3.72%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("D37").Value = "'4.53"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").Value = "'0.0348"
$ws.Range("E39").Value = "  -5.04%  "
$ws.Range("E40").Value = "  -7.13%  "
$ws.Range("D41").Value = "'99.52"
$ws.Range("E41").Value = "  +11.93%  "
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("D43").Value = "'69.39"
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  -7.42%  "
$ws.Range("D46").Value = "'114.41"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("D48").Value = "'76.48"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'8.93"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").Value = "'5.27"
$ws.Range("E50").Value = "  -5.23%  "
$ws.Range("E51").Value = "  -2.29%  "
